# Insert a new data row at row 131 (pushing existing rows 131..159 down to 132..160)
# and populate it with a new weekly price record for "Uva" (Red Globe, Provincia de Limarí).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new entire row above the current row 131, shifting rows 131-159 down to 132-160.
$ws.Range("A131:T131").EntireRow.Insert()

# Populate the newly inserted row 131 with the new record's values.
$ws.Range("A131").Value = 11
$ws.Range("B131").Value = "Vega Monumental Concepción"
$ws.Range("C131").Value = "Bíobío"
$ws.Range("D131").Value = 44722
$ws.Range("E131").Value = 8
$ws.Range("F131").Value = "Fruta"
$ws.Range("G131").Value = 100109
$ws.Range("H131").Value = "Uva"
$ws.Range("I131").Value = 100109001
$ws.Range("J131").Value = "Uva"
$ws.Range("K131").Value = "Red Globe"
$ws.Range("L131").Value = "Primera"
$ws.Range("M131").Value = 280
$ws.Range("N131").Value = 9000
$ws.Range("O131").Value = 10000
$ws.Range("P131").Value = 9536
$ws.Range("Q131").Value = "`$/bandeja 18 kilos"
$ws.Range("R131").Value = "Provincia de Limarí"
$ws.Range("S131").Value = 530
$ws.Range("T131").Value = 18
